$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheets
$ws1.Name = "2015 2016"
$ws2.Name = "2017 2018"

# Remove the duplicate "2017-2018" columns (I:O) that were appended to the
# first sheet, restoring it to only contain the 2015-2016 data (A:H).
$ws1.Range("I1:O1048576").Delete()

# Restore the originally-selected cells on each sheet.
$ws1.Range("J9").Select()
$ws2.Activate()
$ws2.Range("M16").Select()
